$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "首开股份"
$ws.Range("A4").Value = "卧龙电驱"
$ws.Range("A5").Value = "长飞光纤"
$ws.Range("A6").Value = "立讯精密"
$ws.Range("A7").Value = "上海建工"
$ws.Range("A8").Value = "三花智控"
$ws.Range("A9").Value = "华胜天成"
$ws.Range("A10").Value = "欧菲光"
$ws.Range("A11").Value = "工业富联"
$ws.Range("A12").Value = "吉视传媒"
$ws.Range("A13").Value = "金发科技"
$ws.Range("A14").Value = "中际旭创"
$ws.Range("A15").Value = "先导智能"
$ws.Range("A16").Value = "赣锋锂业"
$ws.Range("A17").Value = "景兴纸业"
$ws.Range("A18").Value = "均胜电子"
$ws.Range("A19").Value = "烽火通信"
$ws.Range("A20").Value = "永泰能源"
$ws.Range("A21").Value = "海马汽车"
$ws.Range("B5").Value = "首开股份"
$ws.Range("B6").Value = "三花智控"
$ws.Range("B7").Value = "赣锋锂业"
$ws.Range("B8").Value = "金发科技"
$ws.Range("B9").Value = "工业富联"
$ws.Range("B10").Value = "立讯精密"
$ws.Range("B11").Value = "欧菲光"
$ws.Range("B12").Value = "华胜天成"
$ws.Range("B13").Value = "吉视传媒"
$ws.Range("B14").Value = "云南旅游"
$ws.Range("B15").Value = "凯美特气"
$ws.Range("B16").Value = "和而泰"
$ws.Range("B17").Value = "长飞光纤"
$ws.Range("B18").Value = "永泰能源"
$ws.Range("B19").Value = "山河智能"
$ws.Range("B20").Value = "天普股份"
$ws.Range("B21").Value = "均胜电子"
$ws.Range("C2").Value = "山子高科"
$ws.Range("C3").Value = "首开股份"
$ws.Range("C4").Value = "卧龙电驱"
$ws.Range("C7").Value = "欧菲光"
$ws.Range("C8").Value = "上海建工"
$ws.Range("C9").Value = "吉视传媒"
$ws.Range("C10").Value = "天赐材料"
$ws.Range("C11").Value = "山河智能"
$ws.Range("C12").Value = "先导智能"
$ws.Range("C13").Value = "工业富联"
$ws.Range("C14").Value = "利欧股份"
$ws.Range("C15").Value = "立讯精密"
$ws.Range("C16").Value = "中际旭创"
$ws.Range("C17").Value = "均胜电子"
$ws.Range("C18").Value = "景兴纸业"
$ws.Range("C19").Value = "万通发展"
$ws.Range("C20").Value = "露笑科技"
$ws.Range("C21").Value = "北方稀土"
